# Commit: "Added password only test case"
#
# A new test case, verifyThatAdminCannotLogInWithOnlyPasswordFilled, is
# registered in RUNMANAGER (as priority 3, pushing the previously-3rd test
# down to priority 2 and the previously-2nd test up to priority 4) and two
# matching data rows (chrome + firefox) are appended to DATA.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# ---- DATA sheet: test data for the new case --------------------------
$ws2.Range("A8").Value = "verifyThatAdminCannotLogInWithOnlyPasswordFilled"
$ws2.Range("B8").Value = "yes"
$ws2.Range("C8").Value = "Admin"
$ws2.Range("D8").Value = "admin123"
$ws2.Range("E8").Value = "amuthan"
$ws2.Range("F8").Value = "chrome"

$ws2.Range("A9").Value = "verifyThatAdminCannotLogInWithOnlyPasswordFilled"
$ws2.Range("B9").Value = "yes"
$ws2.Range("C9").Value = "Admin"
$ws2.Range("D9").Value = "admin123"
$ws2.Range("E9").Value = "amuthan"
$ws2.Range("F9").Value = "firefox"

$ws2.Range("A9").Select()

# ---- RUNMANAGER sheet: register the new test case ---------------------
# New row 5: the new test case, at priority 3.
$ws1.Range("A5").Value = "verifyThatAdminCannotLogInWithOnlyPasswordFilled"
$ws1.Range("B5").Value = "To check this test is executed"
$ws1.Range("C5").Value = "yes"
$ws1.Range("D5").Value = "'3"
$ws1.Range("E5").Value = "'1"

# Existing row 4 (verifyThatAdminCannotLogInWithOnlyUsernameFilled) moves
# from priority 3 down to priority 2.
$ws1.Range("D4").Value = "'2"

# Existing row 2 (verifyThatAdminCanLogInWithValidCredentials) moves from
# priority 2 up to priority 4.
$ws1.Range("D2").Value = "'4"

$ws1.Activate()
$ws1.Range("E6").Select()
